$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1480
$ws.Range("I41").Value = 1648.3334
$ws.Range("K41").Value = 1648.3334
$ws.Range("M41").Value = -1208.3334

$ws.Range("H131").Value = 2215.6155
$ws.Range("I131").Value = 2159
$ws.Range("K131").Value = 6477
$ws.Range("M131").Value = -1437

$ws.Range("H135").Value = 1086.6666
$ws.Range("I135").Value = 1061
$ws.Range("J135").Value = 1215
$ws.Range("K135").Value = 9549
$ws.Range("L135").Value = 10935
$ws.Range("M135").Value = -7014
$ws.Range("N135").Value = -16005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3067.6
$ws.Range("I32").Value = 1085.5294
$ws.Range("K32").Value = 1085.5294
$ws.Range("M32").Value = -798.5293999999999

$ws.Range("H74").Value = 4025.7
$ws.Range("I74").Value = 3920
$ws.Range("J74").Value = 4448.5
$ws.Range("K74").Value = 3920
$ws.Range("L74").Value = 4448.5
$ws.Range("M74").Value = -3046
$ws.Range("N74").Value = -6196.5

$ws.Range("H77").Value = 4025.7
$ws.Range("I77").Value = 3920
$ws.Range("J77").Value = 4448.5
$ws.Range("K77").Value = 19600
$ws.Range("L77").Value = 22242.5
$ws.Range("M77").Value = -15232
$ws.Range("N77").Value = -30978.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 2028.1562
$ws.Range("I132").Value = 2031.069
$ws.Range("K132").Value = 6093.207
$ws.Range("M132").Value = -3563.207

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 374.66666
$ws.Range("I22").Value = 361.5
$ws.Range("J22").Value = 401
$ws.Range("K22").Value = 361.5
$ws.Range("L22").Value = 401
$ws.Range("M22").Value = -188.5
$ws.Range("N22").Value = -747

$ws.Range("H36").Value = 659.1429000000001
$ws.Range("I36").Value = 659.1429000000001
$ws.Range("K36").Value = 659.1429000000001
$ws.Range("M36").Value = -125.1429000000001

$ws.Range("H134").Value = 1354.5652
$ws.Range("I134").Value = 1157.091
$ws.Range("K134").Value = 3471.273
$ws.Range("M134").Value = -936.2729999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 762.3333
$ws.Range("I5").Value = 199.8
$ws.Range("J5").Value = 1465.5
$ws.Range("K5").Value = 199.8
$ws.Range("L5").Value = 1465.5
$ws.Range("M5").Value = -87.80000000000001
$ws.Range("N5").Value = -1689.5

$ws.Range("H8").Value = 799
$ws.Range("I8").Value = 799
$ws.Range("K8").Value = 799
$ws.Range("M8").Value = -659

$ws.Range("H11").Value = 187.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 187.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 187.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -467.5

$ws.Range("H12").Value = 1247.5555
$ws.Range("I12").Value = 2065.6
$ws.Range("J12").Value = 225
$ws.Range("K12").Value = 2065.6
$ws.Range("L12").Value = 225
$ws.Range("M12").Value = -1895.6
$ws.Range("N12").Value = -565

$ws.Range("H22").Value = 2311.889
$ws.Range("J22").Value = 3225
$ws.Range("L22").Value = 3225
$ws.Range("N22").Value = -3925

$ws.Range("H70").Value = 75000
$ws.Range("J70").Value = 75000
$ws.Range("L70").Value = 75000
$ws.Range("N70").Value = -75630

$ws.Range("H73").Value = 75000
$ws.Range("J73").Value = 75000
$ws.Range("L73").Value = 75000
$ws.Range("N73").Value = -77184

$ws.Range("H95").Value = 18150
$ws.Range("J95").Value = 18150
$ws.Range("L95").Value = 18150
$ws.Range("N95").Value = -23642

$ws.Range("H96").Value = 11954
$ws.Range("J96").Value = 11954
$ws.Range("L96").Value = 11954
$ws.Range("N96").Value = -17446

$ws.Range("H99").Value = 1744.5
$ws.Range("I99").Value = 1744.5
$ws.Range("K99").Value = 1744.5
$ws.Range("M99").Value = -246.5

$ws.Range("H126").Value = 1744.5
$ws.Range("I126").Value = 1744.5
$ws.Range("K126").Value = 5233.5
$ws.Range("M126").Value = -2763.5

$ws.Range("H134").Value = 8757.875
$ws.Range("I134").Value = 7890.75
$ws.Range("J134").Value = 9625
$ws.Range("K134").Value = 23672.25
$ws.Range("L134").Value = 28875
$ws.Range("M134").Value = -21137.25
$ws.Range("N134").Value = -33945

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4915.8335
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4915.8335
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 14747.5005
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -15313.5005

$ws.Range("H46").Value = 990
$ws.Range("J46").Value = 987.5
$ws.Range("L46").Value = 2962.5
$ws.Range("N46").Value = -3144.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2050.2727
$ws.Range("J102").Value = 2897.5
$ws.Range("L102").Value = 2897.5
$ws.Range("N102").Value = -6141.5

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 1483.0555
$ws.Range("I122").Value = 1380.0667
$ws.Range("K122").Value = 4140.2001
$ws.Range("M122").Value = -1690.2001

$ws.Range("H126").Value = 3445.6924
$ws.Range("J126").Value = 3997.5
$ws.Range("L126").Value = 11992.5
$ws.Range("N126").Value = -16932.5

$ws.Range("H132").Value = 50797.145
$ws.Range("I132").Value = 68895.13
$ws.Range("J132").Value = 5552.1665
$ws.Range("K132").Value = 206685.39
$ws.Range("L132").Value = 16656.4995
$ws.Range("M132").Value = -204155.39
$ws.Range("N132").Value = -21716.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3134.3845
$ws.Range("I7").Value = 2683.5557
$ws.Range("K7").Value = 2683.5557
$ws.Range("M7").Value = -2571.5557

$ws.Range("H22").Value = 957.6667
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -704
$ws.Range("N22").Value = -1465

$ws.Range("H27").Value = 957.6667
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 875
$ws.Range("M27").Value = -892
$ws.Range("N27").Value = -1089

$ws.Range("H40").Value = 4379.8
$ws.Range("I40").Value = 4379.8
$ws.Range("K40").Value = 4379.8
$ws.Range("M40").Value = -4243.8

$ws.Range("H126").Value = 3134.3845
$ws.Range("I126").Value = 2683.5557
$ws.Range("K126").Value = 8050.6671
$ws.Range("M126").Value = -5580.6671

$ws.Range("H132").Value = 5483.45
$ws.Range("I132").Value = 4154.2144
$ws.Range("K132").Value = 12462.6432
$ws.Range("M132").Value = -9932.643199999999

$ws.Range("H136").Value = 3694.6667
$ws.Range("I136").Value = 2084
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 6252
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -3702
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H52").Value = 24013.666
$ws.Range("I52").Value = 13021
$ws.Range("K52").Value = 13021
$ws.Range("M52").Value = -12795

$ws.Range("H122").Value = 6057.143
$ws.Range("I122").Value = 4800
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 14400
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -11950
$ws.Range("N122").Value = -25900

$ws.Range("H126").Value = 3605.7334
$ws.Range("I126").Value = 1509.5555
$ws.Range("K126").Value = 4528.666499999999
$ws.Range("M126").Value = -2058.666499999999

$ws.Range("H132").Value = 1049.3636
$ws.Range("I132").Value = 1049.3636
$ws.Range("K132").Value = 3148.0908
$ws.Range("M132").Value = -618.0907999999999
